$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: previously held the SUM formula - replace with new data entry
$ws.Range("B11").ClearContents()
$ws.Range("A11").Value = 43819
$ws.Range("A11").NumberFormat = "d-mmm"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = "selection operators + writing"

# Row 10: new entry - maken van heuristic crossover + heuristic mutation
$ws.Range("A10").Value = 43815
$ws.Range("A10").NumberFormat = "d-mmm"
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = "maken van heuristic crossover + heuristic mutation"

# Row 17: new SUM total over the full range
$ws.Range("B17").Formula = "=SUM(B2:B16)"

# Update the active selection to match the saved view
$ws.Range("G20").Select() | Out-Null
